$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 14,20
$arr[0,0] = 11
$arr[0,1] = 'Vega Monumental Concepción'
$arr[0,2] = 'Bíobío'
$arr[0,3] = 44902
$arr[0,4] = 8
$arr[0,5] = 'Fruta'
$arr[0,6] = 100103
$arr[0,7] = 'Frutos de hueso (carozo)'
$arr[0,8] = 100103003
$arr[0,9] = 'Damasco'
$arr[0,10] = 'Castle Brite'
$arr[0,11] = 'Primera'
$arr[0,12] = 100
$arr[0,13] = 15000
$arr[0,14] = 16000
$arr[0,15] = 15500
$arr[0,16] = '$/caja 10 kilos'
$arr[0,17] = 'Región de O''Higgins'
$arr[0,18] = 1550
$arr[0,19] = 10
$arr[1,0] = 11
$arr[1,1] = 'Vega Monumental Concepción'
$arr[1,2] = 'Bíobío'
$arr[1,3] = 44902
$arr[1,4] = 8
$arr[1,5] = 'Fruta'
$arr[1,6] = 100103
$arr[1,7] = 'Frutos de hueso (carozo)'
$arr[1,8] = 100103003
$arr[1,9] = 'Damasco'
$arr[1,10] = 'Castle Brite'
$arr[1,11] = 'Segunda'
$arr[1,12] = 50
$arr[1,13] = 13000
$arr[1,14] = 13000
$arr[1,15] = 13000
$arr[1,16] = '$/caja 10 kilos'
$arr[1,17] = 'Región de O''Higgins'
$arr[1,18] = 1300
$arr[1,19] = 10
$arr[2,0] = 11
$arr[2,1] = 'Vega Monumental Concepción'
$arr[2,2] = 'Bíobío'
$arr[2,3] = 44189
$arr[2,4] = 8
$arr[2,5] = 'Fruta'
$arr[2,6] = 100103
$arr[2,7] = 'Frutos de hueso (carozo)'
$arr[2,8] = 100103003
$arr[2,9] = 'Damasco'
$arr[2,10] = 'Dina'
$arr[2,11] = 'Primera'
$arr[2,12] = 200
$arr[2,13] = 15000
$arr[2,14] = 16000
$arr[2,15] = 15500
$arr[2,16] = '$/caja 15 kilos granel'
$arr[2,17] = 'Región de O''Higgins'
$arr[2,18] = 1033
$arr[2,19] = 15
$arr[3,0] = 11
$arr[3,1] = 'Vega Monumental Concepción'
$arr[3,2] = 'Bíobío'
$arr[3,3] = 44189
$arr[3,4] = 8
$arr[3,5] = 'Fruta'
$arr[3,6] = 100103
$arr[3,7] = 'Frutos de hueso (carozo)'
$arr[3,8] = 100103003
$arr[3,9] = 'Damasco'
$arr[3,10] = 'Dina'
$arr[3,11] = 'Segunda'
$arr[3,12] = 100
$arr[3,13] = 14000
$arr[3,14] = 14000
$arr[3,15] = 14000
$arr[3,16] = '$/caja 15 kilos granel'
$arr[3,17] = 'Región de O''Higgins'
$arr[3,18] = 933
$arr[3,19] = 15
$arr[4,0] = 11
$arr[4,1] = 'Vega Monumental Concepción'
$arr[4,2] = 'Bíobío'
$arr[4,3] = 44901
$arr[4,4] = 8
$arr[4,5] = 'Fruta'
$arr[4,6] = 100103
$arr[4,7] = 'Frutos de hueso (carozo)'
$arr[4,8] = 100103003
$arr[4,9] = 'Damasco'
$arr[4,10] = 'Castle Brite'
$arr[4,11] = 'Primera'
$arr[4,12] = 100
$arr[4,13] = 15000
$arr[4,14] = 16000
$arr[4,15] = 15500
$arr[4,16] = '$/caja 10 kilos'
$arr[4,17] = 'Región de O''Higgins'
$arr[4,18] = 1550
$arr[4,19] = 10
$arr[5,0] = 11
$arr[5,1] = 'Vega Monumental Concepción'
$arr[5,2] = 'Bíobío'
$arr[5,3] = 44159
$arr[5,4] = 8
$arr[5,5] = 'Fruta'
$arr[5,6] = 100103
$arr[5,7] = 'Frutos de hueso (carozo)'
$arr[5,8] = 100103003
$arr[5,9] = 'Damasco'
$arr[5,10] = 'Castle Brite'
$arr[5,11] = 'Primera'
$arr[5,12] = 100
$arr[5,13] = 14000
$arr[5,14] = 15000
$arr[5,15] = 14500
$arr[5,16] = '$/caja 15 kilos'
$arr[5,17] = 'Región Metropolitana'
$arr[5,18] = 967
$arr[5,19] = 15
$arr[6,0] = 11
$arr[6,1] = 'Vega Monumental Concepción'
$arr[6,2] = 'Bíobío'
$arr[6,3] = 44187
$arr[6,4] = 8
$arr[6,5] = 'Fruta'
$arr[6,6] = 100103
$arr[6,7] = 'Frutos de hueso (carozo)'
$arr[6,8] = 100103003
$arr[6,9] = 'Damasco'
$arr[6,10] = 'Dina'
$arr[6,11] = 'Primera'
$arr[6,12] = 100
$arr[6,13] = 15000
$arr[6,14] = 16000
$arr[6,15] = 15500
$arr[6,16] = '$/caja 18 kilos'
$arr[6,17] = 'Región Metropolitana'
$arr[6,18] = 861
$arr[6,19] = 18
$arr[7,0] = 11
$arr[7,1] = 'Vega Monumental Concepción'
$arr[7,2] = 'Bíobío'
$arr[7,3] = 44875
$arr[7,4] = 8
$arr[7,5] = 'Fruta'
$arr[7,6] = 100103
$arr[7,7] = 'Frutos de hueso (carozo)'
$arr[7,8] = 100103003
$arr[7,9] = 'Damasco'
$arr[7,10] = 'Castle Brite'
$arr[7,11] = 'Primera'
$arr[7,12] = 50
$arr[7,13] = 31000
$arr[7,14] = 32000
$arr[7,15] = 31400
$arr[7,16] = '$/bandeja 10 kilos'
$arr[7,17] = 'Provincia de Limarí'
$arr[7,18] = 3140
$arr[7,19] = 10
$arr[8,0] = 11
$arr[8,1] = 'Vega Monumental Concepción'
$arr[8,2] = 'Bíobío'
$arr[8,3] = 44579
$arr[8,4] = 8
$arr[8,5] = 'Fruta'
$arr[8,6] = 100103
$arr[8,7] = 'Frutos de hueso (carozo)'
$arr[8,8] = 100103003
$arr[8,9] = 'Damasco'
$arr[8,10] = 'Modesto'
$arr[8,11] = 'Primera'
$arr[8,12] = 180
$arr[8,13] = 13000
$arr[8,14] = 14000
$arr[8,15] = 13444
$arr[8,16] = '$/caja 18 kilos'
$arr[8,17] = 'Región Metropolitana'
$arr[8,18] = 747
$arr[8,19] = 18
$arr[9,0] = 11
$arr[9,1] = 'Vega Monumental Concepción'
$arr[9,2] = 'Bíobío'
$arr[9,3] = 44908
$arr[9,4] = 8
$arr[9,5] = 'Fruta'
$arr[9,6] = 100103
$arr[9,7] = 'Frutos de hueso (carozo)'
$arr[9,8] = 100103003
$arr[9,9] = 'Damasco'
$arr[9,10] = 'Albaricoque'
$arr[9,11] = 'Primera'
$arr[9,12] = 100
$arr[9,13] = 20000
$arr[9,14] = 22000
$arr[9,15] = 21000
$arr[9,16] = '$/caja 18 kilos granel'
$arr[9,17] = 'Región de O''Higgins'
$arr[9,18] = 1167
$arr[9,19] = 18
$arr[10,0] = 11
$arr[10,1] = 'Vega Monumental Concepción'
$arr[10,2] = 'Bíobío'
$arr[10,3] = 44559
$arr[10,4] = 8
$arr[10,5] = 'Fruta'
$arr[10,6] = 100103
$arr[10,7] = 'Frutos de hueso (carozo)'
$arr[10,8] = 100103003
$arr[10,9] = 'Damasco'
$arr[10,10] = 'Modesto'
$arr[10,11] = 'Primera'
$arr[10,12] = 100
$arr[10,13] = 19000
$arr[10,14] = 20000
$arr[10,15] = 19500
$arr[10,16] = '$/caja 18 kilos'
$arr[10,17] = 'Región de O''Higgins'
$arr[10,18] = 1083
$arr[10,19] = 18
$arr[11,0] = 11
$arr[11,1] = 'Vega Monumental Concepción'
$arr[11,2] = 'Bíobío'
$arr[11,3] = 44559
$arr[11,4] = 8
$arr[11,5] = 'Fruta'
$arr[11,6] = 100103
$arr[11,7] = 'Frutos de hueso (carozo)'
$arr[11,8] = 100103003
$arr[11,9] = 'Damasco'
$arr[11,10] = 'Modesto'
$arr[11,11] = 'Segunda'
$arr[11,12] = 50
$arr[11,13] = 18000
$arr[11,14] = 18000
$arr[11,15] = 18000
$arr[11,16] = '$/caja 18 kilos'
$arr[11,17] = 'Región de O''Higgins'
$arr[11,18] = 1000
$arr[11,19] = 18
$arr[12,0] = 11
$arr[12,1] = 'Vega Monumental Concepción'
$arr[12,2] = 'Bíobío'
$arr[12,3] = 44545
$arr[12,4] = 8
$arr[12,5] = 'Fruta'
$arr[12,6] = 100103
$arr[12,7] = 'Frutos de hueso (carozo)'
$arr[12,8] = 100103003
$arr[12,9] = 'Damasco'
$arr[12,10] = 'Castle Brite'
$arr[12,11] = 'Primera'
$arr[12,12] = 100
$arr[12,13] = 18000
$arr[12,14] = 19000
$arr[12,15] = 18500
$arr[12,16] = '$/caja 15 kilos'
$arr[12,17] = 'Región de O''Higgins'
$arr[12,18] = 1233
$arr[12,19] = 15
$arr[13,0] = 11
$arr[13,1] = 'Vega Monumental Concepción'
$arr[13,2] = 'Bíobío'
$arr[13,3] = 44545
$arr[13,4] = 8
$arr[13,5] = 'Fruta'
$arr[13,6] = 100103
$arr[13,7] = 'Frutos de hueso (carozo)'
$arr[13,8] = 100103003
$arr[13,9] = 'Damasco'
$arr[13,10] = 'Castle Brite'
$arr[13,11] = 'Segunda'
$arr[13,12] = 50
$arr[13,13] = 17000
$arr[13,14] = 17000
$arr[13,15] = 17000
$arr[13,16] = '$/caja 15 kilos'
$arr[13,17] = 'Región de O''Higgins'
$arr[13,18] = 1133
$arr[13,19] = 15

$ws.Range("A2:T15").Value = $arr

$ws.Range("D2:D15").NumberFormat = $ws.Range("D2").NumberFormat

Write-Output "Updated data rows 2-15"